$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the mis-computed trial-duration columns (Q..V = fix/mask1/mask2/prime/
# mask3/target duration) for every practice trial row (2-41). All rows share
# the same corrected constants.
$fixDuration    = 1
$mask1Duration  = 0.27
$mask2Duration  = 0.03
$primeDuration  = 0.03
$mask3Duration  = 0.03
$targetDuration = 0.5

for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 17).Value = $fixDuration     # Q
    $ws.Cells.Item($r, 18).Value = $mask1Duration   # R
    $ws.Cells.Item($r, 19).Value = $mask2Duration   # S
    $ws.Cells.Item($r, 20).Value = $primeDuration   # T
    $ws.Cells.Item($r, 21).Value = $mask3Duration   # U
    $ws.Cells.Item($r, 22).Value = $targetDuration  # V
}

# Reflect the author's updated view/selection state on Sheet1.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("S10").Select()
